# Add a new bold "LINK GIT: <url>" line right after the
# "Julia Martins de Almeida Antunes RM98601" paragraph (team-member list),
# matching the commit "Add files via upload".
#
# (The rest of the underlying XML diff is purely Word's interactive
#  spell-checker splitting already-existing sentences into extra runs
#  and wrapping flagged words in <w:proofErr> tags - no visible text
#  changes; that is not something reproducible through legitimate
#  Find/Replace-style COM automation, so it is intentionally left alone.)

$d = $word.ActiveDocument

# Locate the "Julia Martins..." paragraph robustly via Find.
$rng = $d.Content
$found = $rng.Find.Execute("Julia Martins de Almeida Antunes RM98601", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$juliaPara = $rng.Paragraphs(1)
$juliaRange = $juliaPara.Range

# Insert a brand-new empty paragraph right after it.
$juliaRange.InsertParagraphAfter()

# Grab that freshly created paragraph and fill it in.
$newParaRange = $d.Range($juliaRange.End, $juliaRange.End)
$newPara = $newParaRange.Paragraphs(1)

$newPara.Range.Text = "LINK GIT: "
$newPara.Range.Font.Bold = $true
$newPara.Range.LanguageID = "en-US"

# Append the URL right before the paragraph mark, same formatting.
$urlPoint = $d.Range($newPara.Range.End - 1, $newPara.Range.End - 1)
$urlPoint.InsertAfter("https://github.com/JuMartinsDev/GS_Governan-a")
$urlPoint.Font.Bold = $true
$urlPoint.LanguageID = "en-US"
